$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like label into column A without Excel
# auto-converting the text into a real date value/serial number.
function Set-DateLabel($row, $text) {
    $ws.Cells.Item($row, 1).Value = "'" + $text
    $ws.Cells.Item($row, 1).Style = "Normal"
}

# ---- Update existing rows (2025-02-07 / 2025-02-08 data) ----

# Row 27: 2025-02-07 / rel_activity
$ws.Cells.Item(27, 4).Value = 10
$ws.Cells.Item(27, 6).Value = 20

# Row 30: 2025-02-08 / abs_activity
$ws.Cells.Item(30, 3).Value = 10
$ws.Cells.Item(30, 4).Value = 10
$ws.Cells.Item(30, 6).Value = 20

# Row 31: 2025-02-08 / rel_activity
$ws.Cells.Item(31, 4).Value = 10
$ws.Cells.Item(31, 6).Value = 20

# Row 32: 2025-02-08 / abs_sleep
$ws.Cells.Item(32, 4).Value = 9.4
$ws.Cells.Item(32, 6).Value = 9.4

# Row 33: 2025-02-08 / rel_sleep
$ws.Cells.Item(33, 4).Value = 7.184513192693277
$ws.Cells.Item(33, 6).Value = 7.184513192693277

# ---- Append new rows for 2025-02-09 ----

# Row 34: abs_activity
Set-DateLabel 34 "2025-02-09"
$ws.Cells.Item(34, 2).Value = "abs_activity"
$ws.Cells.Item(34, 3).Value = 8.612975372900332
$ws.Cells.Item(34, 4).Value = 9.255128526390363
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 17.86810389929069

# Row 35: rel_activity
Set-DateLabel 35 "2025-02-09"
$ws.Cells.Item(35, 2).Value = "rel_activity"
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0

# Row 36: abs_sleep
Set-DateLabel 36 "2025-02-09"
$ws.Cells.Item(36, 2).Value = "abs_sleep"
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 7.066666666666666
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 7.066666666666666

# Row 37: rel_sleep
Set-DateLabel 37 "2025-02-09"
$ws.Cells.Item(37, 2).Value = "rel_sleep"
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 0
